$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 47, shifting the existing rows 47:53 down to 48:54.
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with this week's price observation.
$ws.Cells.Item(47, 1).Value  = 11
$ws.Cells.Item(47, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(47, 3).Value  = "Bíobío"
$ws.Cells.Item(47, 4).Value  = 44748
$ws.Cells.Item(47, 5).Value  = 8
$ws.Cells.Item(47, 6).Value  = 100112037
$ws.Cells.Item(47, 7).Value  = "Cebollín"
$ws.Cells.Item(47, 8).Value  = "Sin especificar"
$ws.Cells.Item(47, 9).Value  = "Primera"
$ws.Cells.Item(47, 10).Value = 130
$ws.Cells.Item(47, 11).Value = 6000
$ws.Cells.Item(47, 12).Value = 6500
$ws.Cells.Item(47, 13).Value = 6192
$ws.Cells.Item(47, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(47, 15).Value = "Región Metropolitana"
$ws.Cells.Item(47, 16).Value = 172
$ws.Cells.Item(47, 17).Value = 36
$ws.Cells.Item(47, 18).Value = "Hortaliza"
